# Updated the program for Cross Over Design in TA and TE creation
#
# The TE (Trial Elements) sheet's washout row (ET3) was missing the
# washout-period start/end element descriptions. This fills them in and
# fixes the leading space typo on the element name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TE")

# Row 4 corresponds to ETCD = "ET3" (the washout element).
# Correct the element name (drop stray leading space) and populate the
# start/end rule columns that were previously left blank.
$ws.Range("D4").Value = "WASHOUT"
$ws.Range("E4").Value = "End of washout"
$ws.Range("F4").Value = "End of washout period"
$ws.Range("G4").Value = "P7D"
